$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Apply the "0" integer number format to the existing data rows (B2:B21)
$ws.Range("B2:B21").NumberFormat = "0"

# Grow the table by two rows ("MARY GO"/43 and "EL PUENTE"/0)
$row22 = $tbl.ListRows.Add()
$ws.Range("A22").Value = "MARY GO"
$ws.Range("B22").Value = 43
$ws.Range("B22").NumberFormat = "0"

$row23 = $tbl.ListRows.Add()
$ws.Range("A23").Value = "EL PUENTE"
$ws.Range("B23").Value = 0
$ws.Range("B23").NumberFormat = "0"

$ws.Range("D23").Select()
